$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.561.19"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "2.602.56"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.28"
$ws.Range("E5").Value = "  -2.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.39"
$ws.Range("E6").Value = "  +3.42%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +2.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.68"
$ws.Range("E10").Value = "  +2.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.64"
$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").Value = "3.074.99"
$ws.Range("E14").Value = "  -0.45%  "

$ws.Range("D15").Value = "63.428.66"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("E16").Value = "  +7.12%  "

$ws.Range("D17").Value = "2.578.67"
$ws.Range("E17").Value = "  -2.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.51"
$ws.Range("E18").Value = "  +8.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("E19").Value = "  +4.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.31"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.94"
$ws.Range("E21").Value = "  +0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.58"
$ws.Range("E23").Value = "  +2.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.70"
$ws.Range("E24").Value = "  +5.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.34"
$ws.Range("E25").Value = "  +3.19%  "

$ws.Range("E26").Value = "  +0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "568.22"
$ws.Range("E27").Value = "  +3.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("E29").Value = "  +0.46%  "

$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("E31").Value = "  +1.65%  "

$ws.Range("D32").Value = "0.0₃0853"
$ws.Range("E32").Value = "  +1.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.76"
$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.24"
$ws.Range("E34").Value = "  +0.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "167.29"
$ws.Range("E35").Value = "  -1.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.414"
$ws.Range("E36").Value = "  +3.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("E38").Value = "  +3.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.94"
$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "168.69"
$ws.Range("E41").Value = "  +2.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.65"
$ws.Range("E42").Value = "  -0.51%  "

$ws.Range("E43").Value = "  +5.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0588"
$ws.Range("E44").Value = "  +4.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.30"
$ws.Range("E45").Value = "  +2.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.632"
$ws.Range("E46").Value = "  +1.38%  "

$ws.Range("E47").Value = "  +5.78%  "

$ws.Range("E48").Value = "  +4.61%  "

$ws.Range("E49").Value = "  +1.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.25"
$ws.Range("E50").Value = "  +3.53%  "

$ws.Range("D51").Value = "0.0₆0235"
$ws.Range("E51").Value = "  +19.42%  "
